$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prefixes": drop the duplicate/unused prefixes (dc/dcterms duplicate,
# owl, rdf , rdfs ), rename dct -> dcterms, and add the two new geo_core / sf
# prefixes used by the new Traffic / geo_core vocabularies.
# ---------------------------------------------------------------------------
$wsPrefixes = $wb.Worksheets.Item("Prefixes")

$wsPrefixes.Cells.Item(6,1).Value = "dcterms"

$wsPrefixes.Rows.Item(18).Delete()
$wsPrefixes.Rows.Item(17).Delete()
$wsPrefixes.Rows.Item(16).Delete()
$wsPrefixes.Rows.Item(10).Delete()
$wsPrefixes.Rows.Item(9).Delete()

$lastRow = $wsPrefixes.UsedRange.Rows.Count
$wsPrefixes.Cells.Item($lastRow+1,1).Value = "geo_core"
$wsPrefixes.Cells.Item($lastRow+1,2).Value = "https://datos.ign.es/def/geo_core#"
$wsPrefixes.Cells.Item($lastRow+2,1).Value = "sf"
$wsPrefixes.Cells.Item($lastRow+2,2).Value = "http://www.opengis.net/ont/sf#"

$wsPrefixes.Range("A24").Select()

# ---------------------------------------------------------------------------
# Sheet "Source": content unchanged, only the remembered selection moves.
# ---------------------------------------------------------------------------
$wsSource = $wb.Worksheets.Item("Source")
$wsSource.Range("A7").Select()

# ---------------------------------------------------------------------------
# Sheet "Subject": the idPoint row now maps to the sf:Point class (from the
# new OGC "sf" vocabulary) instead of geosparql:Point.
# ---------------------------------------------------------------------------
$wsSubject = $wb.Worksheets.Item("Subject")
$wsSubject.Cells.Item(3,2).Value = "sf:Point"
$wsSubject.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "PredicateObjectMaps": content unchanged, becomes the active sheet.
# ---------------------------------------------------------------------------
$wsPOM = $wb.Worksheets.Item("PredicateObjectMaps")
$wsPOM.Range("A9").Select()
$wsPOM.Activate()

$wb.Save()
